$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row 45 (inherits formatting/style from row 44 above) and fill
# it in with the next log entry: date, description, and duration in hours.
$ws.Rows(45).Insert()

# Write the date through a formula and then collapse it down to a plain
# value, so it ends up stored as literal text ("2012.7.16") instead of
# being auto-recognized and coerced into a date serial number.
$ws.Range("A45").Formula = '="2012.7.16"'
$ws.Range("A45").Copy()
$ws.Range("A45").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("B45").Value = "调整 prop代码，加入加速 减速属性球。"
$ws.Range("D45").Value = 6

# Update the selected cell to reflect the new last row
$ws.Range("C45").Select()
